$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 291.5
$ws.Range("J33").Value = 548
$ws.Range("L33").Value = 548
$ws.Range("N33").Value = -1006
$ws.Range("H40").Value = 6029.522
$ws.Range("I40").Value = 5255.857
$ws.Range("K40").Value = 5255.857
$ws.Range("M40").Value = -5080.857
$ws.Range("H106").Value = 3580.625
$ws.Range("I106").Value = 3635
$ws.Range("K106").Value = 3635
$ws.Range("M106").Value = -3004
$ws.Range("H113").Value = 3018.4
$ws.Range("J113").Value = 3019.6
$ws.Range("L113").Value = 3019.6
$ws.Range("N113").Value = -9527.6
$ws.Range("H126").Value = 77742.75
$ws.Range("J126").Value = 77742.75
$ws.Range("L126").Value = 77742.75
$ws.Range("N126").Value = -87622.75
$ws.Range("H138").Value = 329362.2
$ws.Range("J138").Value = 389535.5
$ws.Range("L138").Value = 1168606.5
$ws.Range("N138").Value = -1178886.5

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1630.4073
$ws.Range("I2").Value = 1501
$ws.Range("K2").Value = 1501
$ws.Range("M2").Value = -1388
$ws.Range("H4").Value = 250
$ws.Range("I4").Value = 250
$ws.Range("K4").Value = 250
$ws.Range("M4").Value = -134
$ws.Range("H32").Value = 12973.013
$ws.Range("I32").Value = 8333.138999999999
$ws.Range("J32").Value = 34515.285
$ws.Range("K32").Value = 8333.138999999999
$ws.Range("L32").Value = 34515.285
$ws.Range("M32").Value = -8046.138999999999
$ws.Range("N32").Value = -35089.285
$ws.Range("H61").Value = 7363.6665
$ws.Range("J61").Value = 0
$ws.Range("L61").Value = 0
$ws.Range("N61").Value = $null
$ws.Range("H63").Value = 2581.6428
$ws.Range("J63").Value = 3899
$ws.Range("L63").Value = 3899
$ws.Range("N63").Value = -5271
$ws.Range("H66").Value = 2581.6428
$ws.Range("J66").Value = 3899
$ws.Range("L66").Value = 19495
$ws.Range("N66").Value = -26359
$ws.Range("H74").Value = 431504.06
$ws.Range("I74").Value = 467046.1
$ws.Range("K74").Value = 467046.1
$ws.Range("M74").Value = -466172.1
$ws.Range("H77").Value = 431504.06
$ws.Range("I77").Value = 467046.1
$ws.Range("K77").Value = 2335230.5
$ws.Range("M77").Value = -2330862.5
$ws.Range("H116").Value = 1630.4073
$ws.Range("I116").Value = 1501
$ws.Range("K116").Value = 1501
$ws.Range("M116").Value = 793
$ws.Range("H132").Value = 2557.3462
$ws.Range("I132").Value = 2009.1765
$ws.Range("J132").Value = 3592.7778
$ws.Range("K132").Value = 6027.529500000001
$ws.Range("L132").Value = 10778.3334
$ws.Range("M132").Value = -3497.529500000001
$ws.Range("N132").Value = -15838.3334
$ws.Range("H136").Value = 7363.6665
$ws.Range("J136").Value = 0
$ws.Range("L136").Value = 0
$ws.Range("N136").Value = $null

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1630.4073
$ws.Range("I3").Value = 1501
$ws.Range("K3").Value = 1501
$ws.Range("M3").Value = -1387
$ws.Range("H82").Value = 89945.25
$ws.Range("I82").Value = 77499
$ws.Range("J82").Value = 94094
$ws.Range("K82").Value = 77499
$ws.Range("L82").Value = 94094
$ws.Range("M82").Value = -77116
$ws.Range("N82").Value = -94860
$ws.Range("H85").Value = 89945.25
$ws.Range("I85").Value = 77499
$ws.Range("J85").Value = 94094
$ws.Range("K85").Value = 77499
$ws.Range("L85").Value = 94094
$ws.Range("M85").Value = -76173
$ws.Range("N85").Value = -96746

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 6516.8184
$ws.Range("I7").Value = 7018.5
$ws.Range("K7").Value = 7018.5
$ws.Range("M7").Value = -6905.5
$ws.Range("H31").Value = 3211400.2
$ws.Range("I31").Value = 5213.82
$ws.Range("J31").Value = 14715952
$ws.Range("K31").Value = 5213.82
$ws.Range("L31").Value = 14715952
$ws.Range("M31").Value = -4918.82
$ws.Range("N31").Value = -14716542
$ws.Range("H34").Value = 3211400.2
$ws.Range("I34").Value = 5213.82
$ws.Range("J34").Value = 14715952
$ws.Range("K34").Value = 5213.82
$ws.Range("L34").Value = 14715952
$ws.Range("M34").Value = -5011.82
$ws.Range("N34").Value = -14716356
$ws.Range("H52").Value = 95847
$ws.Range("J52").Value = 95847
$ws.Range("L52").Value = 95847
$ws.Range("N52").Value = -96435
$ws.Range("H99").Value = 9363.833000000001
$ws.Range("J99").Value = 5888.3076
$ws.Range("L99").Value = 5888.3076
$ws.Range("N99").Value = -8884.3076
$ws.Range("H107").Value = 1777.6
$ws.Range("I107").Value = 1777.6
$ws.Range("K107").Value = 1777.6
$ws.Range("M107").Value = 142.4000000000001
$ws.Range("H126").Value = 9363.833000000001
$ws.Range("J126").Value = 5888.3076
$ws.Range("L126").Value = 17664.9228
$ws.Range("N126").Value = -22604.9228
$ws.Range("H132").Value = 10003401
$ws.Range("I132").Value = 10872218
$ws.Range("J132").Value = 12000
$ws.Range("K132").Value = 32616654
$ws.Range("L132").Value = 36000
$ws.Range("M132").Value = -32614124
$ws.Range("N132").Value = -41060
$ws.Range("H141").Value = 512378.75
$ws.Range("J141").Value = 512378.75
$ws.Range("L141").Value = 512378.75
$ws.Range("N141").Value = -522738.75

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H81").Value = 6500
$ws.Range("J81").Value = 6500
$ws.Range("L81").Value = 19500
$ws.Range("N81").Value = -21746
$ws.Range("H84").Value = 6500
$ws.Range("J84").Value = 6500
$ws.Range("L84").Value = 58500
$ws.Range("N84").Value = -69732
$ws.Range("H122").Value = 1446.7142
$ws.Range("J122").Value = 1544.2778
$ws.Range("L122").Value = 13898.5002
$ws.Range("N122").Value = -18798.5002
$ws.Range("H129").Value = 4806.08
$ws.Range("J129").Value = 3138.4211
$ws.Range("L129").Value = 9415.263300000001
$ws.Range("N129").Value = -19415.2633

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H124").Value = 39779.5
$ws.Range("J124").Value = 39779.5
$ws.Range("L124").Value = 39779.5
$ws.Range("N124").Value = -49599.5
$ws.Range("H132").Value = 2903
$ws.Range("I132").Value = 3210
$ws.Range("K132").Value = 9630
$ws.Range("M132").Value = -7100

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 6140.391
$ws.Range("I7").Value = 4381.4
$ws.Range("K7").Value = 4381.4
$ws.Range("M7").Value = -4269.4
$ws.Range("H126").Value = 6140.391
$ws.Range("I126").Value = 4381.4
$ws.Range("K126").Value = 13144.2
$ws.Range("M126").Value = -10674.2
$ws.Range("H132").Value = 4931.5
$ws.Range("I132").Value = 3929.6667
$ws.Range("K132").Value = 11789.0001
$ws.Range("M132").Value = -9259.000100000001
$ws.Range("H133").Value = 129998
$ws.Range("J133").Value = 129998
$ws.Range("L133").Value = 129998
$ws.Range("N133").Value = -135058
$ws.Range("H136").Value = 8998.5
$ws.Range("I136").Value = 8333
$ws.Range("J136").Value = 10995
$ws.Range("K136").Value = 24999
$ws.Range("L136").Value = 32985
$ws.Range("M136").Value = -22449
$ws.Range("N136").Value = -38085

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 47619892
$ws.Range("I100").Value = 716.9091
$ws.Range("K100").Value = 1433.8182
$ws.Range("M100").Value = -892.8181999999999
$ws.Range("H103").Value = 21000
$ws.Range("J103").Value = 21000
$ws.Range("L103").Value = 21000
$ws.Range("N103").Value = -23344
$ws.Range("H126").Value = 2750.3
$ws.Range("I126").Value = 2567
$ws.Range("K126").Value = 7701
$ws.Range("M126").Value = -5231
$ws.Range("H132").Value = 3834883.2
$ws.Range("I132").Value = 4633105
$ws.Range("J132").Value = 3419.5334
$ws.Range("K132").Value = 13899315
$ws.Range("L132").Value = 10258.6002
$ws.Range("M132").Value = -13896785
$ws.Range("N132").Value = -15318.6002

Write-Host "Applied all Gilgamesh_Profits updates"